$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 2, 3, 6, 7
$ws.Range("F2").Value = -8
$ws.Range("F3").Value = -4
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = -1
